$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1485.9166
$ws.Range("I19").Value = 624
$ws.Range("J19").Value = 2347.8333
$ws.Range("K19").Value = 624
$ws.Range("L19").Value = 2347.8333
$ws.Range("M19").Value = -449
$ws.Range("N19").Value = -2697.8333

$ws.Range("H100").Value = 1875
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -4082

$ws.Range("H107").Value = 954.7857
$ws.Range("I107").Value = 760.0909
$ws.Range("K107").Value = 760.0909
$ws.Range("M107").Value = 1159.9091

$ws.Range("H123").Value = 68000
$ws.Range("J123").Value = 68000
$ws.Range("L123").Value = 68000
$ws.Range("N123").Value = -77800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4315.0713
$ws.Range("I61").Value = 3122.1
$ws.Range("J61").Value = 7297.5
$ws.Range("K61").Value = 3122.1
$ws.Range("L61").Value = 7297.5
$ws.Range("M61").Value = -2910.1
$ws.Range("N61").Value = -7721.5

$ws.Range("H63").Value = 5926.5713
$ws.Range("I63").Value = 5874.5
$ws.Range("J63").Value = 5996
$ws.Range("K63").Value = 5874.5
$ws.Range("L63").Value = 5996
$ws.Range("M63").Value = -5188.5
$ws.Range("N63").Value = -7368

$ws.Range("H66").Value = 5926.5713
$ws.Range("I66").Value = 5874.5
$ws.Range("J66").Value = 5996
$ws.Range("K66").Value = 29372.5
$ws.Range("L66").Value = 29980
$ws.Range("M66").Value = -25940.5
$ws.Range("N66").Value = -36844

$ws.Range("H74").Value = 4296.375
$ws.Range("I74").Value = 4296.375
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4296.375
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3422.375
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 4296.375
$ws.Range("I77").Value = 4296.375
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 21481.875
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -17113.875
$ws.Range("N77").ClearContents()

$ws.Range("H102").Value = 2586.6
$ws.Range("I102").Value = 1911.3334
$ws.Range("J102").Value = 3599.5
$ws.Range("K102").Value = 1911.3334
$ws.Range("L102").Value = 3599.5
$ws.Range("M102").Value = -289.3334
$ws.Range("N102").Value = -6843.5

$ws.Range("H109").Value = 67203.22
$ws.Range("J109").Value = 67203.22
$ws.Range("L109").Value = 67203.22
$ws.Range("N109").Value = -69977.22

$ws.Range("H115").Value = 29600
$ws.Range("J115").Value = 34500
$ws.Range("L115").Value = 34500
$ws.Range("N115").Value = -37634

$ws.Range("H122").Value = 1834.6364
$ws.Range("I122").Value = 1912.6875
$ws.Range("J122").Value = 1626.5
$ws.Range("K122").Value = 5738.0625
$ws.Range("L122").Value = 4879.5
$ws.Range("M122").Value = -3288.0625
$ws.Range("N122").Value = -9779.5

$ws.Range("H132").Value = 2401.8125
$ws.Range("I132").Value = 1522
$ws.Range("J132").Value = 2929.7
$ws.Range("K132").Value = 4566
$ws.Range("L132").Value = 8789.099999999999
$ws.Range("M132").Value = -2036
$ws.Range("N132").Value = -13849.1

$ws.Range("H136").Value = 4315.0713
$ws.Range("I136").Value = 3122.1
$ws.Range("J136").Value = 7297.5
$ws.Range("K136").Value = 9366.299999999999
$ws.Range("L136").Value = 21892.5
$ws.Range("M136").Value = -6816.299999999999
$ws.Range("N136").Value = -26992.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1303.32
$ws.Range("I20").Value = 1245.3
$ws.Range("J20").Value = 1535.4
$ws.Range("K20").Value = 1245.3
$ws.Range("L20").Value = 1535.4
$ws.Range("M20").Value = -998.3
$ws.Range("N20").Value = -2029.4

$ws.Range("H99").Value = 1505.4615
$ws.Range("I99").Value = 1471.75
$ws.Range("K99").Value = 1471.75
$ws.Range("M99").Value = 26.25

$ws.Range("H107").Value = 2116.3076
$ws.Range("I107").Value = 1738.3334
$ws.Range("J107").Value = 2966.75
$ws.Range("K107").Value = 1738.3334
$ws.Range("L107").Value = 2966.75
$ws.Range("M107").Value = 181.6666
$ws.Range("N107").Value = -6806.75

$ws.Range("H122").Value = 68000
$ws.Range("J122").Value = 68000
$ws.Range("L122").Value = 68000
$ws.Range("N122").Value = -77800

$ws.Range("H134").Value = 8241.556
$ws.Range("I134").Value = 8528.208000000001
$ws.Range("J134").Value = 5948.3335
$ws.Range("K134").Value = 25584.624
$ws.Range("L134").Value = 17845.0005
$ws.Range("M134").Value = -23049.624
$ws.Range("N134").Value = -22915.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3939.2307
$ws.Range("I31").Value = 1168.3334
$ws.Range("J31").Value = 6314.2856
$ws.Range("K31").Value = 1168.3334
$ws.Range("L31").Value = 6314.2856
$ws.Range("M31").Value = -873.3334
$ws.Range("N31").Value = -6904.2856

$ws.Range("H34").Value = 3939.2307
$ws.Range("I34").Value = 1168.3334
$ws.Range("J34").Value = 6314.2856
$ws.Range("K34").Value = 1168.3334
$ws.Range("L34").Value = 6314.2856
$ws.Range("M34").Value = -966.3334
$ws.Range("N34").Value = -6718.2856

$ws.Range("H99").Value = 2256.625
$ws.Range("I99").Value = 2110.6
$ws.Range("K99").Value = 2110.6
$ws.Range("M99").Value = -612.5999999999999

$ws.Range("H122").Value = 1353.6428
$ws.Range("I122").Value = 1440.875
$ws.Range("K122").Value = 4322.625
$ws.Range("M122").Value = -1872.625

$ws.Range("H126").Value = 2256.625
$ws.Range("I126").Value = 2110.6
$ws.Range("K126").Value = 6331.799999999999
$ws.Range("M126").Value = -3861.799999999999

$ws.Range("H132").Value = 2627.8333
$ws.Range("I132").Value = 1117.1818
$ws.Range("K132").Value = 3351.5454
$ws.Range("M132").Value = -821.5454

$ws.Range("H134").Value = 2736.0715
$ws.Range("I134").Value = 2283.8333
$ws.Range("K134").Value = 6851.499899999999
$ws.Range("M134").Value = -4316.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13533564
$ws.Range("J131").Value = 24606.066
$ws.Range("L131").Value = 73818.198
$ws.Range("N131").Value = -83898.198

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 12091.625
$ws.Range("J92").Value = 12091.625
$ws.Range("L92").Value = 12091.625
$ws.Range("N92").Value = -15835.625

$ws.Range("H122").Value = 2636.25
$ws.Range("I122").Value = 2272.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6817.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4367.5
$ws.Range("N122").Value = -13900

$ws.Range("H126").Value = 65552.375
$ws.Range("I126").Value = 3294.1538
$ws.Range("K126").Value = 9882.4614
$ws.Range("M126").Value = -7412.4614

$ws.Range("H132").Value = 4196.4116
$ws.Range("I132").Value = 3404.2
$ws.Range("J132").Value = 5328.143
$ws.Range("K132").Value = 10212.6
$ws.Range("L132").Value = 15984.429
$ws.Range("M132").Value = -7682.599999999999
$ws.Range("N132").Value = -21044.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4206.273
$ws.Range("I7").Value = 1974.6923
$ws.Range("K7").Value = 1974.6923
$ws.Range("M7").Value = -1862.6923

$ws.Range("H61").Value = 2693.6428
$ws.Range("J61").Value = 3275.75
$ws.Range("L61").Value = 3275.75
$ws.Range("N61").Value = -3679.75

$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51996

$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159984

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H113").Value = 2693.6428
$ws.Range("J113").Value = 3275.75
$ws.Range("L113").Value = 3275.75
$ws.Range("N113").Value = -7615.75

$ws.Range("H126").Value = 4206.273
$ws.Range("I126").Value = 1974.6923
$ws.Range("K126").Value = 5924.0769
$ws.Range("M126").Value = -3454.0769

$ws.Range("H132").Value = 2823.5
$ws.Range("I132").Value = 1483.5714
$ws.Range("K132").Value = 4450.7142
$ws.Range("M132").Value = -1920.7142

$ws.Range("H136").Value = 3469.6128
$ws.Range("I136").Value = 2285.05
$ws.Range("K136").Value = 6855.150000000001
$ws.Range("M136").Value = -4305.150000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 70062
$ws.Range("J49").Value = 70062
$ws.Range("L49").Value = 70062
$ws.Range("N49").Value = -70522

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H123").Value = 47508.8
$ws.Range("J123").Value = 47508.8
$ws.Range("L123").Value = 47508.8
$ws.Range("N123").Value = -57308.8

$ws.Range("H136").Value = 1996.9362
$ws.Range("I136").Value = 1721.6487
$ws.Range("J136").Value = 3015.5
$ws.Range("K136").Value = 5164.9461
$ws.Range("L136").Value = 9046.5
$ws.Range("M136").Value = -2614.9461
$ws.Range("N136").Value = -14146.5
